$wb = $excel.ActiveWorkbook

$wsCooc = $wb.Worksheets.Item("Cooccurrence")
$wsAssoc = $wb.Worksheets.Item("Associations")

# ---------------------------------------------------------------------------
# Sheet "Cooccurrence": populate header row + 6 data rows (source/target/count)
# ---------------------------------------------------------------------------
$coocHeaders = @("source", "target", "count")
for ($col = 1; $col -le $coocHeaders.Length; $col++) {
    $wsCooc.Cells.Item(1, $col).Value = $coocHeaders[$col - 1]
}

$coocRows = @(
    @("人工智能", "新质生产力", 1),
    @("创新驱动发展", "科技成果转化", 1),
    @("量子计算", "量子通信", 1),
    @("人工智能", "科技成果转化", 1),
    @("人工智能", "元宇宙", 1),
    @("人工智能", "知识产权保护", 1)
)

for ($i = 0; $i -lt $coocRows.Length; $i++) {
    $r = $i + 2
    $row = $coocRows[$i]
    $wsCooc.Cells.Item($r, 1).Value = $row[0]
    $wsCooc.Cells.Item($r, 2).Value = $row[1]
    $wsCooc.Cells.Item($r, 3).Value = $row[2]
}

# Copy the header formatting (bold, border, centered - style index 2) from the
# already-styled "Associations" header row onto the new "Cooccurrence" header.
$wsAssoc.Range("A1:B1").Copy()
$wsCooc.Range("A1:C1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Sheet "Associations": replace placeholder row 2 and append rows 3-15
# ---------------------------------------------------------------------------
$assocRows = @(
    @("新质生产力", 1),
    @("人工智能", 6),
    @("创新驱动发展", 1),
    @("科技成果转化", 4),
    @("量子通信", 1),
    @("量子计算", 2),
    @("知识产权保护", 2),
    @("元宇宙", 2),
    @("氢能", 1),
    @("科技体制改革", 1),
    @("国际人才服务", 1),
    @("碳纤维复合材料", 1),
    @("技术标准互认", 1),
    @("生物降解材料", 1)
)

for ($i = 0; $i -lt $assocRows.Length; $i++) {
    $r = $i + 2
    $row = $assocRows[$i]
    $wsAssoc.Cells.Item($r, 1).Value = $row[0]
    $wsAssoc.Cells.Item($r, 2).Value = $row[1]
}

Write-Host "Cooccurrence and Associations sheets updated."
